# Auto-generated edit script applying the scheduled-runner price/profit updates
# to the Sheets workbook (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 39.5
$ws.Range("I12").Value = 39.5
$ws.Range("K12").Value = 39.5
$ws.Range("M12").Value = 130.5
$ws.Range("H32").Value = 15750
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 29500
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 29500
$ws.Range("M32").Value = -1674
$ws.Range("N32").Value = -30152
$ws.Range("H69").Value = 18315.666
$ws.Range("I69").Value = 18497
$ws.Range("J69").Value = 18225
$ws.Range("K69").Value = 55491
$ws.Range("L69").Value = 54675
$ws.Range("M69").Value = -54617
$ws.Range("N69").Value = -56423
$ws.Range("H72").Value = 18315.666
$ws.Range("I72").Value = 18497
$ws.Range("J72").Value = 18225
$ws.Range("K72").Value = 166473
$ws.Range("L72").Value = 164025
$ws.Range("M72").Value = -162105
$ws.Range("N72").Value = -172761
$ws.Range("H76").Value = 13949.182
$ws.Range("J76").Value = 13393.167
$ws.Range("L76").Value = 13393.167
$ws.Range("N76").Value = -14023.167
$ws.Range("H79").Value = 13949.182
$ws.Range("J79").Value = 13393.167
$ws.Range("L79").Value = 13393.167
$ws.Range("N79").Value = -15577.167
$ws.Range("H98").Value = 2406.2856
$ws.Range("I98").Value = 2250.3845
$ws.Range("K98").Value = 2250.3845
$ws.Range("M98").Value = -752.3845000000001
$ws.Range("H100").Value = 6084.6665
$ws.Range("I100").Value = 6084.6665
$ws.Range("K100").Value = 6084.6665
$ws.Range("M100").Value = -5543.6665
$ws.Range("H107").Value = 1138.7646
$ws.Range("I107").Value = 886.4286
$ws.Range("K107").Value = 886.4286
$ws.Range("M107").Value = 1033.5714
$ws.Range("H112").Value = 3976.8
$ws.Range("I112").Value = 1128.3334
$ws.Range("J112").Value = 5197.5713
$ws.Range("K112").Value = 3385.0002
$ws.Range("L112").Value = 15592.7139
$ws.Range("M112").Value = -2277.0002
$ws.Range("N112").Value = -17808.7139
$ws.Range("H113").Value = 1812.2858
$ws.Range("I113").Value = 1812.2858
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1812.2858
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1441.7142
$ws.Range("N113").ClearContents()
$ws.Range("H115").Value = 9065.166999999999
$ws.Range("I115").Value = 9165
$ws.Range("K115").Value = 27495
$ws.Range("M115").Value = -25928
$ws.Range("H122").Value = 2406.2856
$ws.Range("I122").Value = 2250.3845
$ws.Range("K122").Value = 6751.1535
$ws.Range("M122").Value = -4301.1535
$ws.Range("H129").Value = 1425.5385
$ws.Range("J129").Value = 1040
$ws.Range("L129").Value = 3120
$ws.Range("N129").Value = -13120
$ws.Range("H132").Value = 3314.4
$ws.Range("I132").Value = 3151.4546
$ws.Range("K132").Value = 9454.363799999999
$ws.Range("M132").Value = -6924.363799999999
$ws.Range("H138").Value = 2869.6458
$ws.Range("J138").Value = 3210.6943
$ws.Range("L138").Value = 9632.082900000001
$ws.Range("N138").Value = -19912.0829

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2681.4546
$ws.Range("I2").Value = 2531.4375
$ws.Range("K2").Value = 2531.4375
$ws.Range("M2").Value = -2418.4375
$ws.Range("H32").Value = 7318.0713
$ws.Range("I32").Value = 677.88635
$ws.Range("K32").Value = 677.88635
$ws.Range("M32").Value = -390.88635
$ws.Range("H45").Value = 10136.63
$ws.Range("I45").Value = 12063.904
$ws.Range("K45").Value = 12063.904
$ws.Range("M45").Value = -11686.904
$ws.Range("H88").Value = 3534.1428
$ws.Range("J88").Value = 3807.8
$ws.Range("L88").Value = 3807.8
$ws.Range("N88").Value = -4619.8
$ws.Range("H91").Value = 3534.1428
$ws.Range("J91").Value = 3807.8
$ws.Range("L91").Value = 3807.8
$ws.Range("N91").Value = -6615.8
$ws.Range("H110").Value = 2643.5833
$ws.Range("J110").Value = 2746.5
$ws.Range("L110").Value = 2746.5
$ws.Range("N110").Value = -6836.5
$ws.Range("H116").Value = 2681.4546
$ws.Range("I116").Value = 2531.4375
$ws.Range("K116").Value = 2531.4375
$ws.Range("M116").Value = -237.4375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2681.4546
$ws.Range("I3").Value = 2531.4375
$ws.Range("K3").Value = 2531.4375
$ws.Range("M3").Value = -2417.4375
$ws.Range("H107").Value = 2636.5293
$ws.Range("I107").Value = 2344.4285
$ws.Range("J107").Value = 3999.6667
$ws.Range("K107").Value = 2344.4285
$ws.Range("L107").Value = 3999.6667
$ws.Range("M107").Value = -424.4285
$ws.Range("N107").Value = -7839.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 15097.272
$ws.Range("I16").Value = 2476.3333
$ws.Range("K16").Value = 2476.3333
$ws.Range("M16").Value = -2189.3333
$ws.Range("H107").Value = 1007.2727
$ws.Range("I107").Value = 1006.7
$ws.Range("K107").Value = 1006.7
$ws.Range("M107").Value = 913.3
$ws.Range("H113").Value = 15097.272
$ws.Range("I113").Value = 2476.3333
$ws.Range("K113").Value = 2476.3333
$ws.Range("M113").Value = -306.3332999999998
$ws.Range("H132").Value = 9167.143
$ws.Range("I132").Value = 8390.444
$ws.Range("K132").Value = 25171.332
$ws.Range("M132").Value = -22641.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 69927
$ws.Range("J37").Value = 69927
$ws.Range("L37").Value = 209781
$ws.Range("N37").Value = -210005
$ws.Range("H55").Value = 4174.5713
$ws.Range("J55").Value = 4174.5713
$ws.Range("L55").Value = 12523.7139
$ws.Range("N55").Value = -12877.7139
$ws.Range("H100").Value = 14997.667
$ws.Range("J100").Value = 14997.667
$ws.Range("L100").Value = 44993.001
$ws.Range("N100").Value = -46615.001
$ws.Range("H117").Value = 3303.2
$ws.Range("I117").Value = 1933.3334
$ws.Range("J117").Value = 3890.2856
$ws.Range("K117").Value = 5800.0002
$ws.Range("L117").Value = 11670.8568
$ws.Range("M117").Value = -2358.0002
$ws.Range("N117").Value = -18554.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 32000
$ws.Range("J95").Value = 32000
$ws.Range("L95").Value = 32000
$ws.Range("N95").Value = -37492
$ws.Range("H132").Value = 5738.121
$ws.Range("I132").Value = 5812.963
$ws.Range("K132").Value = 17438.889
$ws.Range("M132").Value = -14908.889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 16039.667
$ws.Range("I7").Value = 19718.572
$ws.Range("J7").Value = 8681.857
$ws.Range("K7").Value = 19718.572
$ws.Range("L7").Value = 8681.857
$ws.Range("M7").Value = -19606.572
$ws.Range("N7").Value = -8905.857
$ws.Range("H55").Value = 336.85715
$ws.Range("I55").Value = 396.4375
$ws.Range("J55").Value = 286.6842
$ws.Range("K55").Value = 396.4375
$ws.Range("L55").Value = 286.6842
$ws.Range("M55").Value = -223.4375
$ws.Range("N55").Value = -632.6841999999999
$ws.Range("H61").Value = 2917.5
$ws.Range("I61").Value = 2821.5386
$ws.Range("K61").Value = 2821.5386
$ws.Range("M61").Value = -2619.5386
$ws.Range("H113").Value = 2917.5
$ws.Range("I113").Value = 2821.5386
$ws.Range("K113").Value = 2821.5386
$ws.Range("M113").Value = -651.5385999999999
$ws.Range("H126").Value = 16039.667
$ws.Range("I126").Value = 19718.572
$ws.Range("J126").Value = 8681.857
$ws.Range("K126").Value = 59155.716
$ws.Range("L126").Value = 26045.571
$ws.Range("M126").Value = -56685.716
$ws.Range("N126").Value = -30985.571
$ws.Range("H136").Value = 1353.2307
$ws.Range("I136").Value = 1099.3636
$ws.Range("K136").Value = 3298.0908
$ws.Range("M136").Value = -748.0907999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H132").Value = 3151.5334
$ws.Range("I132").Value = 2553.05
$ws.Range("K132").Value = 7659.150000000001
$ws.Range("M132").Value = -5129.150000000001
